$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1.75
$ws.Range("H3").Value = 6
$ws.Range("I3").Value = 6.6
$ws.Range("J3").Value = 3.6
$ws.Range("K3").Value = 3.95
